$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the J3 component row (original row 5) - rows below shift up by one
$ws.Rows.Item(5).Delete()

# Re-assert clean Mid X / Mid Y / Rotation values for every data row so that
# no floating point drift is introduced by the row shift above.
$ws.Range("B2").Value = 135
$ws.Range("C2").Value = -90.73
$ws.Range("E2").Value = 90

$ws.Range("B3").Value = 138
$ws.Range("C3").Value = -99.46
$ws.Range("E3").Value = 180

$ws.Range("B4").Value = 155.88
$ws.Range("C4").Value = -117.38
$ws.Range("E4").Value = 0

$ws.Range("B5").Value = 127.8475
$ws.Range("C5").Value = -93
$ws.Range("E5").Value = 180

$ws.Range("B6").Value = 144.71
$ws.Range("C6").Value = -106.6
$ws.Range("E6").Value = 180

$ws.Range("B7").Value = 140.91999999999999
$ws.Range("C7").Value = -107.5
$ws.Range("E7").Value = -90

$ws.Range("B8").Value = 149.0275
$ws.Range("C8").Value = -106.753
$ws.Range("E8").Value = 90

$ws.Range("B9").Value = 147.45750000000001
$ws.Range("C9").Value = -106.74299999999999
$ws.Range("E9").Value = 90

$ws.Range("B10").Value = 157.77000000000001
$ws.Range("C10").Value = -112.88
$ws.Range("E10").Value = -90

$ws.Range("B11").Value = 155.94999999999999
$ws.Range("C11").Value = -112.91
$ws.Range("E11").Value = -90

$ws.Range("B12").Value = 140.88
$ws.Range("C12").Value = -104.33
$ws.Range("E12").Value = 180

$ws.Range("B13").Value = 140.85
$ws.Range("C13").Value = -89.84
$ws.Range("E13").Value = 0

# Update the current selection to match the saved view state
$ws.Range("C10").Select()
